# Generate Report for Handoff
#
# The localization XLIFF handoff files finished generating for both
# target languages, so the per-language "Status" moves from
# "In Translation" to "Ready for handoff", and the associated
# generate/handoff timestamps advance a few seconds. The "Status"
# columns are widened slightly so the new, longer status text fits.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Status: "In Translation" -> "Ready for handoff" ---
$wsOverview.Range("E2").Value = "Ready for handoff"   # zh-cn status
$wsOverview.Range("F2").Value = "Ready for handoff"   # de-de status
$wsZhCn.Range("C2").Value = "Ready for handoff"       # Status
$wsDeDe.Range("C2").Value = "Ready for handoff"       # Status

# --- Timestamps advance as the handoff XLIFFs are (re)generated ---
$wsOverview.Range("G2").Value = "2016-09-04 08:42:52" # Latest HO Xliff Generate Date
$wsDeDe.Range("H2").Value = "2016-09-04 08:42:52"     # Latest Handoff Datetime (de-de)
$wsZhCn.Range("H2").Value = "2016-09-04 08:42:47"     # Latest Handoff Datetime (zh-cn)

# --- Widen the Status columns to fit the new, longer status text ---
$wsOverview.Columns.Item(5).ColumnWidth = 16.33  # zh-cn status column
$wsOverview.Columns.Item(6).ColumnWidth = 16.33  # de-de status column
$wsZhCn.Columns.Item(3).ColumnWidth = 16.33       # Status column
$wsDeDe.Columns.Item(3).ColumnWidth = 16.33       # Status column
